$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 1.471908333333333
$ws.Range("N2").Value = 4.415725
$ws.Range("O2").Value = 0.2507683239665115
$ws.Range("P2").Value = 0.2507683239665115
$ws.Range("Q2").Value = 296.5008358955889
$ws.Range("R2").Value = 2668.5075230603
$ws.Range("S2").Value = 0.1212088752342817
$ws.Range("T2").Value = 0.1212088752342817
# Row 3
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("N3").Value = 7.040756999999999
$ws.Range("O3").Value = 0.3998434758381655
$ws.Range("P3").Value = 0.3998434758381655
$ws.Range("Q3").Value = 472.7627594195106
$ws.Range("R3").Value = 4254.864834775596
$ws.Range("S3").Value = 0.1932643533661846
$ws.Range("T3").Value = 0.1932643533661847
# Row 4
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.8927443333333333
$ws.Range("N4").Value = 2.678233
$ws.Range("O4").Value = 0.152096428242656
$ws.Range("P4").Value = 0.1520964282426559
$ws.Range("Q4").Value = 179.8341887737916
$ws.Range("R4").Value = 1618.507698964124
$ws.Range("S4").Value = 0.07351581213624853
$ws.Range("T4").Value = 0.07351581213624853
# Row 5
$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("M5").Value = 1.158022666666667
$ws.Range("N5").Value = 3.474068
$ws.Range("O5").Value = 0.1972917719526671
$ws.Range("P5").Value = 0.197291771952667
$ws.Range("Q5").Value = 233.2717879755005
$ws.Range("R5").Value = 2099.446091779504
$ws.Range("S5").Value = 0.09536098257192437
$ws.Range("T5").Value = 0.09536098257192437
# Row 6
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("M6").Value = 1.471908333333333
$ws.Range("N6").Value = 4.415725
$ws.Range("O6").Value = 0.2507683239665115
$ws.Range("P6").Value = 0.2507683239665115
$ws.Range("Q6").Value = 96.28836419757224
$ws.Range("R6").Value = 866.59527777815
$ws.Range("S6").Value = 0.03936246684527557
$ws.Range("T6").Value = 0.03936246684527558
# Row 7
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("N7").Value = 7.040756999999999
$ws.Range("O7").Value = 0.3998434758381655
$ws.Range("P7").Value = 0.3998434758381655
$ws.Range("S7").Value = 0.06276241477404998
$ws.Range("T7").Value = 0.06276241477404999
# Row 8
$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.8927443333333333
$ws.Range("N8").Value = 2.678233
$ws.Range("O8").Value = 0.152096428242656
$ws.Range("P8").Value = 0.1520964282426559
$ws.Range("Q8").Value = 58.40098160776689
$ws.Range("R8").Value = 525.6088344699019
$ws.Range("S8").Value = 0.02387419000649337
$ws.Range("T8").Value = 0.02387419000649337
# Row 9
$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("M9").Value = 1.158022666666667
$ws.Range("N9").Value = 3.474068
$ws.Range("O9").Value = 0.1972917719526671
$ws.Range("P9").Value = 0.197291771952667
$ws.Range("Q9").Value = 75.75479107759911
$ws.Range("R9").Value = 681.7931196983919
$ws.Range("S9").Value = 0.03096838830956022
$ws.Range("T9").Value = 0.03096838830956022
# Row 10
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 1.471908333333333
$ws.Range("N10").Value = 4.415725
$ws.Range("O10").Value = 0.2507683239665115
$ws.Range("P10").Value = 0.2507683239665115
$ws.Range("Q10").Value = 88.95454658221112
$ws.Range("R10").Value = 800.5909192399001
$ws.Range("S10").Value = 0.03636441868920119
$ws.Range("T10").Value = 0.0363644186892012
# Row 11
$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("N11").Value = 7.040756999999999
$ws.Range("O11").Value = 0.3998434758381655
$ws.Range("P11").Value = 0.3998434758381655
$ws.Range("Q11").Value = 141.8356773871853
$ws.Range("R11").Value = 1276.521096484668
$ws.Range("S11").Value = 0.05798210609513139
$ws.Range("T11").Value = 0.05798210609513141
# Row 12
$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.8927443333333333
$ws.Range("N12").Value = 2.678233
$ws.Range("O12").Value = 0.152096428242656
$ws.Range("P12").Value = 0.1520964282426559
$ws.Range("Q12").Value = 53.95286213623245
$ws.Range("R12").Value = 485.575759226092
$ws.Range("S12").Value = 0.02205580876509189
$ws.Range("T12").Value = 0.02205580876509189
# Row 13
$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("M13").Value = 1.158022666666667
$ws.Range("N13").Value = 3.474068
$ws.Range("O13").Value = 0.1972917719526671
$ws.Range("P13").Value = 0.197291771952667
$ws.Range("Q13").Value = 69.98491612040355
$ws.Range("R13").Value = 629.864245083632
$ws.Range("S13").Value = 0.02860967639668589
$ws.Range("T13").Value = 0.02860967639668589
# Row 14
$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("M14").Value = 1.471908333333333
$ws.Range("N14").Value = 4.415725
$ws.Range("O14").Value = 0.2507683239665115
$ws.Range("P14").Value = 0.2507683239665115
$ws.Range("Q14").Value = 131.6850763253472
$ws.Range("R14").Value = 1185.165686928125
$ws.Range("S14").Value = 0.05383256319775299
$ws.Range("T14").Value = 0.053832563197753
# Row 15
$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("N15").Value = 7.040756999999999
$ws.Range("O15").Value = 0.3998434758381655
$ws.Range("P15").Value = 0.3998434758381655
$ws.Range("Q15").Value = 209.9683795827917
$ws.Range("R15").Value = 1889.715416245125
$ws.Range("S15").Value = 0.08583460160279949
$ws.Range("T15").Value = 0.0858346016027995
# Row 16
$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.8927443333333333
$ws.Range("N16").Value = 2.678233
$ws.Range("O16").Value = 0.152096428242656
$ws.Range("P16").Value = 0.1520964282426559
$ws.Range("Q16").Value = 79.86985535151389
$ws.Range("R16").Value = 718.8286981636249
$ws.Range("S16").Value = 0.03265061733482216
$ws.Range("T16").Value = 0.03265061733482216
# Row 17
$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("M17").Value = 1.158022666666667
$ws.Range("N17").Value = 3.474068
$ws.Range("O17").Value = 0.1972917719526671
$ws.Range("P17").Value = 0.197291771952667
$ws.Range("Q17").Value = 103.6031251356111
$ws.Range("R17").Value = 932.4281262205001
$ws.Range("S17").Value = 0.04235272467449657
$ws.Range("T17").Value = 0.04235272467449657
